$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("A6").Value = "currency"
$ws.Range("B6").Value = "US Dollar"
$ws.Range("C5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("C12").ClearContents()
